# "Add content examples from old templates, add checklist"
#
# Note: the Swate Excel add-in's own bookkeeping parts (xl/webextensions/*,
# customXml/item1.xml's <SwateTable> "checklist" entry, the workbook
# fileVersion/revisionPtr GUIDs, and the orphaned legacy-comment font in
# styles.xml) are internal artifacts of the Swate task-pane / Excel session
# that produced the commit. They are not reachable through the documented
# Excel object model (CustomXMLParts is unimplemented here and webextension
# parts aren't exposed at all), so this script sticks to the
# user-observable workbook content the diff actually represents.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Expand the annotation table to cover the new example rows (A1:DL6)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:DL6")) | Out-Null

# Hide the Term Source REF / Term Accession Number helper columns for the
# first annotated parameter column, matching the rest of the table.
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(4).Hidden = $true

# Add content examples (rows 2-6) taken from older SWATE templates.
$ws.Range("B2").Value = 'remove leaflets form shoot'
$ws.Range("E2").Value = 'shoot without leaves'
$ws.Range("H2").Value = 'protein localization'
$ws.Range("N2").Value = 'compound treatment'
$ws.Range("Q2").Value = 'Nup107-aa33-51 antibody'
$ws.Range("T2").Value = 'DAPI'
$ws.Range("W2").Value = 'DAPI:nuclei;TRITC:HA_Flag tagged protein'
$ws.Range("Z2").Value = 'ENSG00000117399'
$ws.Range("AC2").Value = 'RuBisCo'
$ws.Range("AF2").Value = 'GRCh37, Ensembl release 61, Feb 2011'
$ws.Range("AI2").Value = 'Sigma-Aldrich, microscope plain, size 25 mm × 75 mm'
$ws.Range("AL2").Value = 'BRAND® L × W 18 mm × 18 mm, square;  Sigma-Aldrich'
$ws.Range("AO2").Value = 'water'
$ws.Range("AR2").Value = '50% v/v ethanol, 3.7% v/v formaldehyde, 5% v/v acetic acid'
$ws.Range("AU2").Value = 'nailpolish'
$ws.Range("AX2").Value = 'gene deletion screen'
$ws.Range("BA2").Value = 'primary screen'
$ws.Range("BD2").Value = 'haploid deletion library'
$ws.Range("BG2").Value = 'Ambion'
$ws.Range("BJ2").Value = 'Bioneer haploid deletion library v.2 modified to generate a GFP-tubulin expressing library (Dixon et al., 2008).'
$ws.Range("BM2").Value = 'JL_120731_S6A'
$ws.Range("BP2").Value = 'A1'
$ws.Range("BS2").Value = 1
$ws.Range("BV2").Value = 's2748'
$ws.Range("BY2").Value = 'CGAAAUGACUAUUACCUGATT'
$ws.Range("CB2").Value = 'UCAGGUAAUAGUCAUUUCGGA'
$ws.Range("CE2").Value = 'GRCh37, Ensembl release 61, Feb 2011'
$ws.Range("CH2").Value = 'DL-alpha-Methyl-p-tyrosine'
$ws.Range("CK2").Value = 78875
$ws.Range("CN2").Value = 'Neurotransmission'
$ws.Range("CQ2").Value = '10 nanogram per milliliter'
$ws.Range("CT2").Value = '1 hour'
$ws.Range("CW2").Value = 9994
$ws.Range("CZ2").Value = 'empty well'
$ws.Range("DC2").Value = 'checks for transfection'
$ws.Range("DF2").Value = 'pass'
$ws.Range("B3").Value = 'extract archargonia from gametophyte'
$ws.Range("H3").Value = 'high content screen'
$ws.Range("N3").Value = 'antibody target'
$ws.Range("Q3").Value = 'Alexa Fluor 546 conjugated secondary antibody'
$ws.Range("W3").Value = 'GFP:endogenous alpha tubulin 2;Cascade blue:growth media'
$ws.Range("Z3").Value = 84240
$ws.Range("AX3").Value = 'compound screen'
$ws.Range("BA3").Value = 'secondary screen'
$ws.Range("BD3").Value = 'compound library'
$ws.Range("BG3").Value = 'Sigma'
$ws.Range("BJ3").Value = 'LOPAC-1280 compound library'
$ws.Range("BM3").Value = 1921
$ws.Range("BP3").Value = 'C3'
$ws.Range("BS3").Value = 'PSORS1C3'
$ws.Range("BV3").Value = 's20068'
$ws.Range("BY3").Value = 'CGGAAGCAGUUCCAACUUUtt'
$ws.Range("CB3").Value = 'AAAGUUGGAACUGCUUCCGtt'
$ws.Range("CH3").Value = 'N-Acetyl-L-Cysteine'
$ws.Range("CK3").Value = 78955
$ws.Range("CN3").Value = 'Glutamate'
$ws.Range("CW3").Value = 168
$ws.Range("CZ3").Value = 'positive control'
$ws.Range("DC3").Value = 'gives strong phenotype'
$ws.Range("DF3").Value = 'fail'
$ws.Range("H4").Value = 'time-lapse imaging'
$ws.Range("AX4").Value = 'RNAi screen'
$ws.Range("BA4").Value = 'validation screen'
$ws.Range("BD4").Value = 'siRNA library'
$ws.Range("CN4").Value = '6-Methoxy-1,2,3,4-tetrahydro-9H-pyrido[3,4b] indole'
$ws.Range("CZ4").Value = 'negative control'
$ws.Range("DC4").Value = 'wild type'
$ws.Range("DF4").Value = 'no cells'
$ws.Range("H5").Value = 'image cytometry'
$ws.Range("AX5").Value = 'protein screen'
$ws.Range("BD5").Value = 'HA-Flag protein fusion library'
$ws.Range("DC5").Value = 'non-targeting siRNA'
$ws.Range("DC6").Value = 'no treatment'

# Restore full-sheet selection artifact left on BioImageArchive_Imaging,
# then return focus to the main annotation sheet.
$wsBia = $wb.Worksheets.Item("BioImageArchive_Imaging")
$wsBia.Activate() | Out-Null
$wsBia.Cells.Select() | Out-Null
$ws.Activate() | Out-Null
